# Updates cryptos list values (Price / Volume(1h) columns, plus two row
# re-rankings) to match the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '35.212.48'
$c.ClearFormats()
$ws.Range('E2').Value = '  +1.14%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.860.31'
$c.ClearFormats()
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  +0.72%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '239.09'
$c.ClearFormats()
$ws.Range('E5').Value = '  +3.58%  '
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  +0.64%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '42.32'
$c.ClearFormats()
$ws.Range('E8').Value = '  +6.80%  '
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('E10').Value = '  +1.43%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0991'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.57%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '2.128.62'
$c.ClearFormats()
$ws.Range('E12').Value = '  +1.12%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '11.48'
$c.ClearFormats()
$ws.Range('E13').Value = '  +1.30%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.677'
$c.ClearFormats()
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '1.839.70'
$c.ClearFormats()
$ws.Range('E15').Value = '  -0.21%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '4.72'
$c.ClearFormats()
$ws.Range('E16').Value = '  +1.53%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '35.177.35'
$c.ClearFormats()
$ws.Range('E17').Value = '  +0.95%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '69.94'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  +1.19%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '241.45'
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.30%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '12.23'
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('E24').Value = '  +0.33%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '169.39'
$c.ClearFormats()
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('E26').Value = '  +24.43%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '8.06'
$c.ClearFormats()
$ws.Range('E27').Value = '  +3.58%  '
$ws.Range('E29').Value = '  +0.09%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.0562'
$c.ClearFormats()
$ws.Range('E30').Value = '  +1.70%  '
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').Value = '  +1.88%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.82'
$c.ClearFormats()
$ws.Range('E33').Value = '  +27.63%  '
$ws.Range('E34').Value = '  +2.03%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '2.06'
$c.ClearFormats()
$ws.Range('E35').Value = '  +9.88%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.817'
$c.ClearFormats()
$ws.Range('E36').Value = '  +17.77%  '
$ws.Range('E37').Value = '  +8.33%  '
$ws.Range('E38').Value = '  +4.19%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.0202'
$c.ClearFormats()
$ws.Range('E39').Value = '  +4.39%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '90.09'
$c.ClearFormats()
$ws.Range('E40').Value = '  -1.33%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.348.96'
$c.ClearFormats()
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '15.08'
$c.ClearFormats()
$ws.Range('E42').Value = '  +3.46%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.0597'
$c.ClearFormats()
$ws.Range('E43').Value = '  +14.42%  '
$ws.Range('E44').Value = '  +1.95%  '
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('B46').Value = 'Gas'
$ws.Range('C46').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '12.31'
$c.ClearFormats()
$ws.Range('E46').Value = '  +47.86%  '
$ws.Range('B47').Value = 'MXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.74'
$c.ClearFormats()
$ws.Range('E47').Value = '  -0.62%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '6.57'
$c.ClearFormats()
$ws.Range('E48').Value = '  +4.77%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.043.20'
$c.ClearFormats()
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0686'
$c.ClearFormats()
$ws.Range('E50').Value = '  +0.73%  '
